$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Ativacao:" date from 01/01/2018 to 01/01/2022 (cells B8/C8).
#    Force text format first so Excel doesn't auto-convert the "dd/mm/yyyy"
#    looking string into a real date serial number.
$dateCells = $ws.Range("B8:C8")
$dateCells.NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = "01/01/2022"
$ws.Cells.Item(8, 3).Value = "01/01/2022"
$dateCells.NumberFormat = "General"

# 2. Update "Objetivos:" text (cells B10/C10)
$ws.Cells.Item(10, 2).Value = "Apresentar aos alunos a Engenharia Bioquímica, as características da profissão e orientar quanto as atribuições e as áreas de atuação do Engenheiro Bioquímico. Além disso, desenvolver nos alunos uma visão macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atuação do Engenheiro Bioquímico na indústria, pesquisa e ensino, e empreendedorismo e inovação em engenharia."
$ws.Cells.Item(10, 3).Value = "Apresentar aos alunos a Engenharia Bioquímica, as características da profissão e orientar quanto as atribuições e as áreas de atuação do Engenheiro Bioquímico. Além disso, desenvolver nos alunos uma visão macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atuação do Engenheiro Bioquímico na indústria, pesquisa e ensino, e empreendedorismo e inovação em engenharia."

# 3. Update "Objectives:" (English) text (cells B11/C11)
$ws.Cells.Item(11, 2).Value = "To present to the Biochemical Engineering student the characteristics of the profession and to guide in relation to the attributes and the action areas of the biochemical engineering. Besides, to develop in the students a macro view of types and stages of an industrial bioprocess and, finally, to guide about the action of the biochemical engineering on the industry, research and teaching, and entrepreneurship and innovation in engineering."
$ws.Cells.Item(11, 3).Value = "To present to the Biochemical Engineering student the characteristics of the profession and to guide in relation to the attributes and the action areas of the biochemical engineering. Besides, to develop in the students a macro view of types and stages of an industrial bioprocess and, finally, to guide about the action of the biochemical engineering on the industry, research and teaching, and entrepreneurship and innovation in engineering."

# 4. Insert a new row after the existing "Docentes responsaveis" teacher row (row 13),
#    to add a second teacher: Valdeir Arantes. This pushes everything from the old
#    row 14 onward down by one row, and copies formatting from row 13.
$ws.Rows("14").Insert()
$ws.Cells.Item(14, 2).Value = "5817181 - Valdeir Arantes"
$ws.Cells.Item(14, 3).Value = "5817181 - Valdeir Arantes"

# 5. Update "Programa:" text (now at row 17, was row 16 before insert) - cells B17/C17
$ws.Cells.Item(17, 2).Value = "1.Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2.Mercado de trabalho da Engenharia do Brasil 3.Atribuições e áreas de atuação do Engenheiro Bioquímico 4.Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5.A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6.Escalas de produção – laboratório, piloto, industrial. 7.Estudo de casos (processos biotecnológicos). 8.Empreendedorismo e Inovação em Engenharia.9.Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso."
$ws.Cells.Item(17, 3).Value = "1.Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2.Mercado de trabalho da Engenharia do Brasil 3.Atribuições e áreas de atuação do Engenheiro Bioquímico 4.Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5.A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6.Escalas de produção – laboratório, piloto, industrial. 7.Estudo de casos (processos biotecnológicos). 8.Empreendedorismo e Inovação em Engenharia.9.Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso."

# 6. Update "Syllabus:" (English) text (now at row 18) - cells B18/C18
$ws.Cells.Item(18, 2).Value = "1.History of the Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2.Job market of Engineering in Brazil3.Attributes and action areas of biochemical engineering4.Definitions and concepts – enzymatic process, general fermentative process, transformation agents, bioreactor, raw material, types of substrates, conversion of substrate into product, types of biotechnological products, products recovery, between others.5.The Bioprocesses Industry – types of industries, equipment, installations, main unit operations6.Production scales – laboratory, pilot, industrial.7.Studies of cases (biotechnological processes).8.Entrepreneurship and Innovation in Engineering.9.Supervised visitation – visits to laboratories and bioprocess industry"
$ws.Cells.Item(18, 3).Value = "1.History of the Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2.Job market of Engineering in Brazil3.Attributes and action areas of biochemical engineering4.Definitions and concepts – enzymatic process, general fermentative process, transformation agents, bioreactor, raw material, types of substrates, conversion of substrate into product, types of biotechnological products, products recovery, between others.5.The Bioprocesses Industry – types of industries, equipment, installations, main unit operations6.Production scales – laboratory, pilot, industrial.7.Studies of cases (biotechnological processes).8.Entrepreneurship and Innovation in Engineering.9.Supervised visitation – visits to laboratories and bioprocess industry"

# 7. Update "Metodo:" text (now at row 20) - cells B20/C20
$ws.Cells.Item(20, 2).Value = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina; exercícios; dinâmicas. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a aplicações dos conceitos abordados à um processo, produto ou serviço na área de Engenharia de Bioquímica e que relacione com a formação acadêmica e atribuições profissionais do Engenheiro Bioquímico."
$ws.Cells.Item(20, 3).Value = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina; exercícios; dinâmicas. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a aplicações dos conceitos abordados à um processo, produto ou serviço na área de Engenharia de Bioquímica e que relacione com a formação acadêmica e atribuições profissionais do Engenheiro Bioquímico."

# 8. Update "Criterio:" text (now at row 21) - cells B21/C21
$ws.Cells.Item(21, 2).Value = "A nota (N) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros."
$ws.Cells.Item(21, 3).Value = "A nota (N) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros."

# 9. Update "Norma de recuperacao:" text (now at row 22) - cells B22/C22
$ws.Cells.Item(22, 2).Value = "Média Final = (N + Prova Recuperação)/2"
$ws.Cells.Item(22, 3).Value = "Média Final = (N + Prova Recuperação)/2"
